$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newQuery = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (diag:diagnosis)-->(c)`nOPTIONAL MATCH (f)-[*]->(samp:sample)`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`nWHERE demo.breed IN ['Parson Russell Terrier'] `n`nOPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)`nWITH`n        f, parent, c, demo, diag, s, samp,`n        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,`n        toInteger(floor(log(f.file_size)/log(1024))) as i,`n        2 as precision`nWITH`n        f, parent, c, demo, diag, s, samp,`n        f.file_size /(1024^i) AS value, `n        10^precision AS factor,`n        units[i] as unit`nWITH    `n        f, parent, c, demo, diag, s, samp, unit,`n        round(factor * value)/factor AS size`nRETURN `n        coalesce(f.file_name, '') AS ``File Name``,`n       coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_type, '') AS ``File Type``,`n       CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,`n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(samp.sample_id, '') AS ``Sample ID``,`n        coalesce(c.case_id, '') AS ``Case ID``,`n        coalesce(demo.breed,'') AS Breed ,`n        coalesce(diag.disease_term,'') AS Diagnosis"

$ws.Range("B4").Value = $newQuery
$ws.Rows.Item(4).RowHeight = 409.5

$ws.Range("C4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
